# Decrement the "剩余" (remaining) column E by 1 for each data row (rows 2-99),
# except row 36 which is left unchanged per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)  # Column E is the 5th column
    $cell.Value2 = $cell.Value2 - 1
}
